$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (R1:T1) ---
# Copy style from the existing last header cell (Q1) so new headers match formatting
$ws.Range("Q1").Copy($ws.Range("R1"))
$ws.Range("Q1").Copy($ws.Range("S1"))
$ws.Range("Q1").Copy($ws.Range("T1"))
$ws.Range("R1").Value = "31/12/2023"
$ws.Range("S1").Value = "31/03/2024"
$ws.Range("T1").Value = "30/06/2024"

# --- Data rows (R2:T80) ---
# $null marks rows that are section headers with no numeric data (columns stay blank,
# matching existing blank cells in B:Q of the same rows)
$rstData = @(
    @(70140837.888, 75030626.30400001, 77114253.31200001),
    @(22160576.512, 24859822.08, 25227261.952),
    @(1345474.944, 1644980.992, 2126797.056),
    @(9630015.488, 10857413.632, 9890046.976),
    @(5106679.808, 5565636.096, 5915739.136),
    @(2890509.056, 2893513.984, 2993146.112),
    @(0, 0, 0),
    @(1354840.96, 1369675.008, 1410204.032),
    @(127725, 393756.992, 415560),
    @(1705329.024, 2134845.056, 2475769.088),
    @(4599224.832, 4761602.048, 5137283.072),
    @(0, 0, 0),
    @(1051288, 967097.9840000001, 1341964.032),
    @(0, 0, 0),
    @(494196.992, 525982.0159999999, 530228),
    @(0, 0, 0),
    @(0, 0, 0),
    @(2090211.968, 2329889.024, 2318185.984),
    @(0, 0, 0),
    @(942, 942, 942),
    @(38223, 36645, 35592),
    @(39826161.664, 41692246.016, 42806902.784),
    @(3516655.104, 3680310.016, 3907212.032),
    @(0, 0, 0),
    @(70140837.888, 75030626.30400001, 77114253.31200001),
    @(15328766.976, 16866872.32, 18124396.544),
    @(727766.976, 762332.992, 787822.0159999999),
    @(6884024.832, 6735143.936, 6351306.24),
    @(500799.008, 527348, 509030.016),
    @(4182392.064, 4848754.176, 6531112.96),
    @(511, 125, 0),
    @(123657, 133710, 4123),
    @(2909616.088, 3859458.896, 3941002.12),
    @(0, 0, 0),
    @(0, 0, 0),
    @(48060338.176, 51441008.64, 52202426.368),
    @(41120432.128, 44414074.88, 45317513.216),
    @(528, 528, 528),
    @(4678538.24, 4640870.912, 4473525.248),
    @(1494779.008, 1614904.064, 1679814.016),
    @(0, 0, 0),
    @(766065.9840000001, 770627.008, 731046.0159999999),
    @(0, 0, 0),
    @(0, 0, 0),
    @(3425486.08, 3505858.048, 3600121.088),
    @(3326246.144, 3216884.736, 3187307.776),
    @(1174361.984, 1174361.984, 1174361.984),
    @(2101220.992, 2101220.992, 2101084.032),
    @(0, 0, 0),
    @(430539.008, 417483.008, 466783.008),
    @(0, 0, 0),
    @(157924, 101865, 66284),
    @(0, 0, 0),
    @(-537800, -578046.0159999999, -621204.992),
    @(0, 0, 0),
    @($null, $null, $null),
    @($null, $null, $null),
    @(8521823.232000001, 9146986.495999999, 10307617.792),
    @(-6884329.472, -6677400.064, -7689915.904),
    @(1637494.272, 2469585.92, 2617701.888),
    @(0, 0, 0),
    @(-1056482.112, -801201.9840000001, -869460.992),
    @(-31300.992, -53232, -131575),
    @(39482.016, -36142, 131044),
    @(0, 0, 0),
    @(2403, -1047, 161),
    @(-1541279.232, -1476937.984, -1509481.984),
    @(299503.936, 356775.008, 363992.992),
    @(-1840782.848, -1833713.024, -1873474.944),
    @($null, $null, $null),
    @($null, $null, $null),
    @($null, $null, $null),
    @(-949682.944, 101025, 238388.992),
    @(359, -70729, -47216),
    @(231296.992, 59806, -32328),
    @($null, $null, $null),
    @($null, $null, $null),
    @(340803.008, -103156, -109547),
    @(-377224, -13054, 49298)
)

for ($i = 0; $i -lt $rstData.Length; $i++) {
    $rowNum = $i + 2
    $vals = $rstData[$i]
    for ($j = 0; $j -lt 3; $j++) {
        $colNum = 18 + $j
        $val = $vals[$j]
        if ($null -eq $val) {
            # Touch the cell so it becomes a real (blank) cell in the sheet,
            # extending the row span / sheet dimension without adding a style or value.
            $ws.Cells.Item($rowNum, $colNum).Font.Bold = $false
        } else {
            $ws.Cells.Item($rowNum, $colNum).Value = $val
        }
    }
}
